# Update risk-driver calibration estimates with values re-computed using a
# std moving-average approach (per commit message).

$wb = $excel.ActiveWorkbook

# --- Linear sheet -----------------------------------------------------
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.005356548018998067
$wsLinear.Range("B3").Value = -0.03363850548963302
$wsLinear.Range("B4").Value = 1.672186344313421
$wsLinear.Range("B5").Value = "[1.0, 0.23778558207983122, 0.022827615227641325, -0.03249033695689084, -0.03176266677585463, 0.04283006738441504, 0.2217164567998397, 0.40089721551891677, 0.20831746620862254, 0.02212732223966535, -0.03230612987206453, -0.04368897576701962, 0.02448350071653799, 0.20592583476159718, 0.3763317198525951, 0.19603237457372225, 0.014629507570053353, -0.026712652893962432, -0.030348085630541775, 0.005185430974261977]"

# --- NonLinear sheet ----------------------------------------------------
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.004987678675128903
$wsNonLinear.Range("B5").Value = -0.04638124571656335
$wsNonLinear.Range("B6").Value = 1.801532830624212
$wsNonLinear.Range("B7").Value = -0.01921237485867761
$wsNonLinear.Range("B8").Value = 0.01511618490025547
$wsNonLinear.Range("B9").Value = 1.557856378060666
$wsNonLinear.Range("B10").Value = "[1.0, 0.23790921385639446, 0.024007684664335307, -0.03140850667671699, -0.03075087992639604, 0.04343299711018512, 0.2214904677698268, 0.39997421599243993, 0.20852048681586038, 0.02293892453716549, -0.031309047737878536, -0.04280330667300946, 0.025376986449402877, 0.20589813093902565, 0.3758023698480617, 0.19585907286182647, 0.015234373561518154, -0.025843726531520506, -0.029260796907132812, 0.0060210359635150855]"
